# Product requirements document.docx - PRD dates edit, actualize README (refs #12)
#
# 1. One of the PRD timeline rows had its date corrected from
#    29.05.2025 to 29.04.2025 (the surrounding "Редактирование и
#    модерация постов..." milestone is actually an April item).
# 2. The "TableGrid" table style's display name was actualized to its
#    localized (Russian) name "Сетка таблицы" used by the current Word UI.

$d = $word.ActiveDocument

# --- 1) Fix the mis-typed milestone date in the second table -------------
$null = $d.Content.Find.Execute(
    "29.05.2025",   # FindText
    $true,          # MatchCase
    $false,         # MatchWholeWord
    $false,         # MatchWildcards
    $false,         # MatchSoundsLike
    $false,         # MatchAllWordForms
    $true,          # Forward
    1,              # Wrap (wdFindContinue)
    $false,         # Format
    "29.04.2025",   # ReplaceWith
    2               # Replace (wdReplaceAll)
)

# --- 2) Actualize the TableGrid style's display name ----------------------
$tableGridStyle = $d.Styles("TableGrid")
$tableGridStyle.NameLocal = "Сетка таблицы"
